$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update cell values (row 2)
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("L2").Value = 0.0125
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 3.18
$ws.Range("Q2").Value = 1

# Update cell values (row 3)
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("M3").Value = 0.0125
$ws.Range("O3").Value = 6.37
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 2

# Update cell values (row 4)
$ws.Range("O4").Value = 9.55
$ws.Range("P4").Value = 2

# Update sheet view: selection (top-left cell scroll position is not
# representable through this COM-interop surface, so only the selection
# is reproduced here).
$ws.Activate()
$ws.Range("Q5").Select()
